# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" > "Impact" bullet list to use
# short, impact-focused accomplishment statements instead of the longer
# job-duty-style bullets (some of which duplicated bullets used elsewhere
# in the resume). Six bullets collapse down to four.

$d = $word.ActiveDocument

# Replace the text of a paragraph (identified by its current index) while
# leaving its paragraph mark - and therefore its paragraph formatting -
# untouched.
function Set-ParagraphText {
    param([int]$Index, [string]$NewText)

    $para = $d.Paragraphs.Item($Index)
    $full = $para.Range
    # Exclude the trailing paragraph-mark character so the replacement
    # stays inside this paragraph instead of swallowing the next one.
    $body = $d.Range($full.Start, $full.End - 1)
    $body.Text = $NewText
}

# Locate the six existing "Impact" bullets under "KEY ACHIEVEMENTS AND
# IMPACT" by their current, distinctive wording. (Two of these bullet
# strings are duplicated verbatim earlier in the Professional Experience
# section, so we scope everything to this heading's own paragraphs rather
# than doing a blind document-wide Find/Replace.)
$startIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($t -eq "Impact") {
        $styleName = $d.Paragraphs.Item($i).Range.ParagraphStyle.NameLocal
        if ($styleName -like "*Heading 3*") {
            $next = $d.Paragraphs.Item($i + 1).Range.Text
            if ($next -like "*Delivered*additional revenue*") {
                $startIndex = $i + 1
                break
            }
        }
    }
}

if ($startIndex -eq -1) {
    throw "Could not locate the Key Achievements 'Impact' bullet list"
}

# startIndex..startIndex+5 are the six bullets being replaced:
#   0: Delivered $4.9M additional revenue ...
#   1: Built redistricting platform ... serving 12,847 analysts ...
#   2: Achieved 87% prediction accuracy ...
#   3: Trigonometric algorithm for boundary estimation ...
#   4: Discovered systematic race coding errors ...
#   5: Developed longitudinal data analysis methods ...
#
# New content keeps the first four bullet slots (rewritten) and removes
# the last two paragraphs entirely.

Set-ParagraphText ($startIndex)     "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
Set-ParagraphText ($startIndex + 1) "• Real-time collaboration at national scale"
Set-ParagraphText ($startIndex + 2) "• Revenue generation: Delivered `$4.9M additional revenue through optimization"
Set-ParagraphText ($startIndex + 3) "• 23% conversion rate improvement"

# Delete the two trailing bullets outright (bottom-up so indices stay
# valid while deleting).
$d.Paragraphs.Item($startIndex + 5).Range.Delete()
$d.Paragraphs.Item($startIndex + 4).Range.Delete()
